$d = $word.ActiveDocument

# Locate the 7-paragraph ".then(...)" promise-chain block that follows the
# "      var alertOptions = { height: 120, width: 260 };" line inside the
# first sample (the "button.onclick" handler), and collapse it down to a
# single statement call plus a stray trailing ";" paragraph, per the diff:
#
#   context.navigation.openAlertDialog(alertStrings, alertOptions).then(
#       function (success) { console.log("Alert dialog closed"); },
#       function (error) { console.log(error.message); });
#
# becomes:
#
#   context.navigation.openAlertDialog(alertStrings, alertOptions);
#   ;
#
# (note: the *second* occurrence of this "...).then(" block, further down
# in the userSettings/userId sample, is left untouched.)

$rng = $d.Content
$rng.Find.Execute("context.navigation.openAlertDialog(alertStrings, alertOptions).then(", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

# Resolve the paragraph index containing the found match.
$startIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $cand = $d.Paragraphs.Item($i)
  if ($cand.Range.Start -le $rng.Start -and $cand.Range.End -ge $rng.End) {
    $startIdx = $i
    break
  }
}

$p1 = $d.Paragraphs.Item($startIdx)
$p7 = $d.Paragraphs.Item($startIdx + 6)

$blockRange = $d.Range($p1.Range.Start, $p7.Range.End)
$blockRange.Delete()

# $startIdx now holds the paragraph that used to follow the deleted block
# (the "  }" closing brace). Insert the two replacement paragraphs
# immediately before it, inheriting its bold+italic run formatting.
$afterPara = $d.Paragraphs.Item($startIdx)
$afterPara.Range.InsertBefore("      context.navigation.openAlertDialog(alertStrings, alertOptions);`r;`r")
